$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cryptos.xlsx refresh -------------------------------------------------
# Updated Price / Volume(1h) figures for the latest data pull, including the
# Dogecoin/Toncoin (rows 9-10) and InjectiveProtocol/ONDO (rows 45-46) rank
# swaps (name + link + price + volume all moved together).
#
# Each target value is entered with a leading apostrophe so the engine keeps
# it as literal text (matching the original inline-string cells) instead of
# silently reinterpreting a number-looking string (e.g. "1.00", "608.65") as
# a numeric value -- exactly how typing an apostrophe-prefixed value works
# in real Excel. The immediately-following Style reset clears the resulting
# quote-prefix flag so the cell keeps its original (default) formatting.

$ws.Range("D2").Value = "'66.481.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.53%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.591.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.77%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'608.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.35%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'149.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.79%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.11%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'Dogecoin"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.136"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.29%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'Toncoin"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'8.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.38%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.78%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.208.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.98%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.60%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'29.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.61%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.627.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.85%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'66.580.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.56%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +0.99%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'11.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.79%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'15.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.74%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'427.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.48%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.619"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.18%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'78.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.38%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.03%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +2.26%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +4.49%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.37%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.62%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.06%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.591.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.89%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.51%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'25.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.67%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.85%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'7.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.16%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.01%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.58%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.82%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'177.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.75%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0857"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.84%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.61%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.896"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.20%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -2.42%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +9.54%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.11%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'InjectiveProtocol"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'25.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.84%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'ONDO"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.89%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'23.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.74%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.20%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.953"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.98%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.429.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.17%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.94%  "
$ws.Range("E51").Style = "Normal"
